# Deploy preview for PR 10 🛫
#
# Replace the plain "first chapter" intro paragraph with a tracked
# insertion (w:ins) of slightly revised wording, authored as "PR Preview".

$d = $word.ActiveDocument

# Tracked changes are attributed to Application.UserName, so set it before
# recording anything.
$word.UserName = "PR Preview"

$oldText = "This is the first chapter of your book. Replace this content with your own."
$q = [char]0x201C
$newText = "This ${q}is${q} the first chapter of your book. Replace this content with your own."

# Locate the original (untracked) run and remove it without recording a
# revision mark - the target markup has no <w:del>, only a fresh <w:ins>.
$d.TrackRevisions = $false
$target = $d.Content
$found = $target.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $found) {
    throw "Could not locate the original first-chapter paragraph text."
}

# Now record the replacement text as a tracked insertion at that same spot.
$d.TrackRevisions = $true
$target.InsertAfter($newText)

$d.TrackRevisions = $false

Write-Output "Inserted tracked-change replacement paragraph for Chapter 1 intro."
